# "names with td prefix" - prefix the English product/module names on the
# tdRPA architecture slide with "td" (tdSelector, tdLocator, tdWorker, tdPower).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeWithText {
    param($slide, [string]$needle)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*$needle*") {
                return $shp
            }
        }
    }
    return $null
}

# --- "Selector - UI" -> "tdSelector - UI" -------------------------------
$shp = Get-ShapeWithText $s "Selector - UI"
$tr = $shp.TextFrame.TextRange
$found = $tr.Find("Selector - UI", 0)
$start = $found.Start
$found.Characters(1, 9).Delete()                          # drop "Selector "
$tr.Characters($start, 0).InsertBefore(" ") | Out-Null     # re-add separating space as its own run
$tr.Characters($start, 0).InsertBefore("tdSelector") | Out-Null

# --- "Locator - UI" -> "tdLocator - UI" ----------------------------------
$shp = Get-ShapeWithText $s "Locator - UI"
$tr = $shp.TextFrame.TextRange
$found = $tr.Find("Locator - UI", 0)
$start = $found.Start
$found.Characters(1, 8).Delete()                           # drop "Locator "
$tr.Characters($start, 0).InsertBefore(" ") | Out-Null
$tr.Characters($start, 0).InsertBefore("tdLocator") | Out-Null

# --- "Worker" -> "tdWorker" ----------------------------------------------
$shp = Get-ShapeWithText $s "Worker"
$tr = $shp.TextFrame.TextRange
$found = $tr.Find("Worker", 0)
$found.Text = "tdWorker"

# --- "Power" -> "tdPower" -------------------------------------------------
$shp = Get-ShapeWithText $s "Power"
$tr = $shp.TextFrame.TextRange
$found = $tr.Find("Power", 0)
$found.Text = "tdPower"
